$wb = $excel.ActiveWorkbook

# OFF sheet: update Home row (row 2) stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 214
$wsOff.Range("C2").Value = 162
$wsOff.Range("D2").Value = 38
$wsOff.Range("E2").Value = 18

# DEF sheet: update Home row (row 2) stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 173
$wsDef.Range("C2").Value = 123
$wsDef.Range("D2").Value = 39
$wsDef.Range("E2").Value = 21
$wsDef.Range("F2").Value = 3
